# Generate Report for Handback
#
# Fills in the "b4286b45-1845-47ff-a03f-494807f1de7a" row (row 7) on both
# the "zh-cn" and "de-de" localization-status sheets now that a (stale)
# handback file showed up for that item: a Latest Target File link, the
# handback xliff name, the handback datetime, and an "out of date" error
# message - plus the hyperlink that goes with the new Latest Target File
# cell.

$wb = $excel.ActiveWorkbook

$errorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/9b774210d0a6add5867b0e1ab7da649ec28a6346/e2e/b4286b45-1845-47ff-a03f-494807f1de7a.md, latest: https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3ce28d2a9c4860184968ddf86df935dd0ffc21a/e2e/b4286b45-1845-47ff-a03f-494807f1de7a.md."
$handoffUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d3ce28d2a9c4860184968ddf86df935dd0ffc21a/e2e/b4286b45-1845-47ff-a03f-494807f1de7a.md"
$handoffDisplay = "b4286b45-1845-47ff-a03f-494807f1de7a.md"

# --- zh-cn sheet, row 7 -----------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I7"), $handoffUrl, [Type]::Missing, [Type]::Missing, $handoffDisplay) | Out-Null
$wsZhCn.Range("I7").Style = "HyperLink"

$wsZhCn.Range("J7").Value = "b4286b45-1845-47ff-a03f-494807f1de7a.0d6f8d7abcef784d676838c987548d7e711241d4.zh-cn.xlf"
$wsZhCn.Range("K7").Value = "2016-08-31 00:58:56"
$wsZhCn.Range("P7").Value = $errorMessage

# --- de-de sheet, row 7 -------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I7"), $handoffUrl, [Type]::Missing, [Type]::Missing, $handoffDisplay) | Out-Null
$wsDeDe.Range("I7").Style = "HyperLink"

$wsDeDe.Range("J7").Value = "b4286b45-1845-47ff-a03f-494807f1de7a.0d6f8d7abcef784d676838c987548d7e711241d4.de-de.xlf"
$wsDeDe.Range("K7").Value = "2016-08-31 00:59:09"
$wsDeDe.Range("P7").Value = $errorMessage

Write-Output "Report generated for handback row 7 (zh-cn, de-de)"
